$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$t = $nm.Theme
$cs = $t.ThemeColorScheme
$c = $cs.Colors(3)
Write-Host "via Colors(3):" $c.RGB
$c.RGB = 111111
Write-Host "after:" $c.RGB
